$wb = $excel.ActiveWorkbook

$wsJindRaw = $wb.Worksheets.Item("JIND_raw")
$wsJindRaw.Range("B2").Value = 1.0
$wsJindRaw.Range("I2").Value = 0.0297029702970297
$wsJindRaw.Range("C3").Value = 0.9741379310344828
$wsJindRaw.Range("C4").Value = 0.013546798029556651
$wsJindRaw.Range("E4").Value = 0.025906735751295335
$wsJindRaw.Range("C5").Value = 0.0012315270935960591
$wsJindRaw.Range("E5").Value = 0.9740932642487047
$wsJindRaw.Range("B6").Value = 0.0
$wsJindRaw.Range("C6").Value = 0.007389162561576354
$wsJindRaw.Range("I9").Value = 0.9603960396039604

$wsJind = $wb.Worksheets.Item("JIND")
$wsJind.Range("B2").Value = 0.9771689497716894
$wsJind.Range("F2").Value = 0.07755102040816327
$wsJind.Range("C3").Value = 0.9605911330049262
$wsJind.Range("D4").Value = 0.9441964285714286
$wsJind.Range("D5").Value = 0.0
$wsJind.Range("D6").Value = 0.008928571428571428
$wsJind.Range("C8").Value = 0.0012315270935960591
$wsJind.Range("I8").Value = 0.9504950495049505
$wsJind.Range("B9").Value = 0.0228310502283105
$wsJind.Range("C9").Value = 0.022167487684729065
$wsJind.Range("F9").Value = 0.04081632653061224
$wsJind.Range("I9").Value = 0.0297029702970297
